$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" (column D) cell's value while forcing it to stay a
# text string (these columns hold values like "64.195.74" that must not be
# coerced into numbers by Excel's normal type inference), and without
# leaving a permanent number-format override on the cell.
function Set-PriceText {
    param(
        [string]$Ref,
        [string]$Value
    )
    $ws.Range($Ref).NumberFormat = "@"
    $ws.Range($Ref).Value = $Value
    $ws.Range($Ref).Style = "Normal"
}

# Helper to set a plain text/percent cell's value directly by A1 reference
function Set-CellValue {
    param(
        [string]$Ref,
        [string]$Value
    )
    $ws.Range($Ref).Value = $Value
}

# Row 2 - Bitcoin
Set-PriceText "D2" "64.239.47"
Set-CellValue "E2" "  -0.87%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.425.73"
Set-CellValue "E3" "  -0.08%  "

# Row 5 - BNB
Set-PriceText "D5" "572.69"
Set-CellValue "E5" "  -0.23%  "

# Row 6 - Solana
Set-PriceText "D6" "161.77"
Set-CellValue "E6" "  +1.53%  "

# Row 7 - USDC
Set-CellValue "E7" "  +0.03%  "

# Row 8 - LidoStakedEther
Set-PriceText "D8" "3.428.58"
Set-CellValue "E8" "  -0.06%  "

# Row 9 - XRP
Set-CellValue "E9" "  -5.57%  "

# Row 10 - Toncoin
Set-PriceText "D10" "7.30"
Set-CellValue "E10" "  +1.56%  "

# Row 11 - Dogecoin
Set-CellValue "E11" "  -1.96%  "

# Row 12 - Cardano
Set-PriceText "D12" "0.425"
Set-CellValue "E12" "  -3.30%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-PriceText "D13" "4.018.57"
Set-CellValue "E13" "  +0.05%  "

# Row 14 - TRON
Set-CellValue "E14" "  +1.23%  "

# Row 15 - Avalanche
Set-PriceText "D15" "27.07"
Set-CellValue "E15" "  -2.37%  "

# Row 16 - ShibaInu
Set-CellValue "E16" "  -6.85%  "

# Row 17 - WrappedBTC
Set-PriceText "D17" "64.285.24"
Set-CellValue "E17" "  -0.88%  "

# Row 18 - WrappedEther
Set-PriceText "D18" "3.500.37"
Set-CellValue "E18" "  +2.67%  "

# Row 19 - Polkadot
Set-PriceText "D19" "6.10"
Set-CellValue "E19" "  -3.98%  "

# Row 20 - Chainlink
Set-PriceText "D20" "13.61"
Set-CellValue "E20" "  -1.82%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "377.36"
Set-CellValue "E21" "  -0.99%  "

# Row 22 - Uniswap
Set-PriceText "D22" "7.84"

# Row 23 - Dai
Set-CellValue "E23" "  -0.22%  "

# Row 24 - Litecoin
Set-PriceText "D24" "71.28"
Set-CellValue "E24" "  -0.93%  "

# Row 25 - Polygon
Set-CellValue "E25" "  -5.20%  "

# Row 26 - PEPE
Set-PriceText "D26" "0.0000117"
Set-CellValue "E26" "  -2.25%  "

# Row 27 - InternetComputer(DFINITY)
Set-CellValue "E27" "  -3.96%  "

# Row 29 - Binance-PegBSC-USD
Set-PriceText "D29" "0.995"
Set-CellValue "E29" "  -0.40%  "

# Row 30 - NEARProtocol
Set-PriceText "D30" "6.03"
Set-CellValue "E30" "  -1.96%  "

# Row 31 - Fetch.AI
Set-PriceText "D31" "1.41"
Set-CellValue "E31" "  -3.66%  "

# Row 32 - PancakeSwap
Set-PriceText "D32" "2.02"
Set-CellValue "E32" "  +0.34%  "

# Row 33 - EthereumClassic
Set-PriceText "D33" "22.96"

# Row 34 - Aptos
Set-PriceText "D34" "7.13"
Set-CellValue "E34" "  +1.48%  "

# Row 35 - ImmutableX
Set-CellValue "E35" "  -4.02%  "

# Row 36 - Monero
Set-PriceText "D36" "159.92"
Set-CellValue "E36" "  -0.66%  "

# Row 37 - Mantle
Set-PriceText "D37" "0.859"
Set-CellValue "E37" "  +11.23%  "

# Row 38 - Stacks
Set-PriceText "D38" "1.82"
Set-CellValue "E38" "  -4.41%  "

# Row 39 - Maker
Set-PriceText "D39" "2.806.77"
Set-CellValue "E39" "  -3.00%  "

# Row 40 - Hedera
Set-PriceText "D40" "0.0729"
Set-CellValue "E40" "  -3.07%  "

# Row 41 - EnergySwap
Set-PriceText "D41" "25.82"
Set-CellValue "E41" "  -2.39%  "

# Row 42 - OKB
Set-CellValue "E42" "  -0.47%  "

# Row 43 - RenderToken
Set-PriceText "D43" "6.49"
Set-CellValue "E43" "  -4.30%  "

# Row 44 - InjectiveProtocol
Set-PriceText "D44" "26.10"
Set-CellValue "E44" "  +0.65%  "

# Row 45 - Filecoin
Set-PriceText "D45" "4.43"
Set-CellValue "E45" "  -3.02%  "

# Row 46 - VeChain
Set-PriceText "D46" "0.0307"
Set-CellValue "E46" "  -2.82%  "

# Row 47 - dogwifhat
Set-PriceText "D47" "2.47"
Set-CellValue "E47" "  +8.30%  "

# Row 48 - Bittensor
Set-PriceText "D48" "336.53"
Set-CellValue "E48" "  +5.83%  "

# Row 49 - ONDO
Set-PriceText "D49" "1.06"
Set-CellValue "E49" "  -2.14%  "

# Row 50 - was Cosmos, now Stellar (rows 50 and 51 swapped)
Set-CellValue "B50" "Stellar"
Set-CellValue "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText "D50" "0.103"
Set-CellValue "E50" "  -2.76%  "

# Row 51 - was Stellar, now Cosmos
Set-CellValue "B51" "Cosmos"
Set-CellValue "C51" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-PriceText "D51" "6.32"
Set-CellValue "E51" "  -3.10%  "
